$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 36
# from serial date 45673 (2025-01-16) to 45674 (2025-01-17).
for ($row = 2; $row -le 36; $row++) {
    $ws.Cells.Item($row, 3).Value = 45674
}
